# Handles float input without breaking stuff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12) ---
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value2 = 14
$ws.Range("C10").Value2 = 3
$ws.Range("D10").Value2 = 11
$ws.Range("E10").Value2 = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value2 = 4
$ws.Range("C11").Value2 = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value2 = 56
$ws.Range("C12").Value2 = -3
$ws.Range("E12").Value2 = "53/112"

# --- Remove the third question block (columns G:H) entirely ---
$ws.Range("G15:H21").Clear()

# --- Remove the per-question answer detail for the second question
#     (columns D:E) everywhere except the first few rows that remain ---
$ws.Range("D19:E40").Clear()

# --- Row 16: second question's "Student Ans" cell now shows "Option C" ---
$ws.Range("D16").Style = "incorrectStyle"
$ws.Range("D16").Value2 = "Option C"

# --- Row 18: mark first question's student answer correct; keep second
#     question's student answer ---
$ws.Range("A18").Style = "correctStyle"
$ws.Range("A18").Value2 = "Option B"
$ws.Range("D18").Style = "correctStyle"
$ws.Range("D18").Value2 = "Option D"

# --- Rows 19-40: first question's "Student Ans" column (A) gets filled in
#     with the student's answers, colored green (correct) or red (incorrect)
#     to match the "Correct Ans" column (B) ---
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A19").Value2 = "Option C"

$ws.Range("A20").Style = "correctStyle"
$ws.Range("A20").Value2 = "Option B"

$ws.Range("A21").Style = "correctStyle"
$ws.Range("A21").Value2 = "Option C"

$ws.Range("A25").Style = "correctStyle"
$ws.Range("A25").Value2 = "Option A"

$ws.Range("A26").Style = "correctStyle"
$ws.Range("A26").Value2 = "Option C"

$ws.Range("A29").Style = "correctStyle"
$ws.Range("A29").Value2 = "Option D"

$ws.Range("A30").Style = "incorrectStyle"
$ws.Range("A30").Value2 = "Option D"

$ws.Range("A31").Style = "correctStyle"
$ws.Range("A31").Value2 = "Option D"

$ws.Range("A33").Style = "correctStyle"
$ws.Range("A33").Value2 = "Option D"

$ws.Range("A34").Style = "correctStyle"
$ws.Range("A34").Value2 = "Option B"

$ws.Range("A36").Style = "incorrectStyle"
$ws.Range("A36").Value2 = "Option D"

$ws.Range("A38").Style = "correctStyle"
$ws.Range("A38").Value2 = "Option A"

$ws.Range("A39").Style = "correctStyle"
$ws.Range("A39").Value2 = "Option D"

$ws.Range("A40").Style = "correctStyle"
$ws.Range("A40").Value2 = "Option D"
